$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ma phong" (room code) header label and its numbering column data.
# The room-code numbering (1..8) in column B and its header text are no longer needed.
$ws.Range("B2").ClearContents()
$ws.Range("B3:B10").ClearContents()

# Highlight the loaded data region (everything except the outer label row/column)
# with the new background fill used to show the subject data has been loaded.
$ws.Range("B2:E20").Interior.Color = 16777214

# Remove the frozen header row/column split.
$excel.ActiveWindow.FreezePanes = $false
